$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 346, shifting existing data (rows 346-439) down to 349-442.
$ws.Rows("346:348").Insert()

# Populate the 3 newly-inserted rows with a new price-report date group
# (2023-10-19, serial 45218) for "Cultivar IV Región" at "Provincia del Elquí".

# Row 346: Especial
$ws.Cells.Item(346, 1).Value = 3
$ws.Cells.Item(346, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(346, 3).Value = "Coquimbo"
$ws.Cells.Item(346, 4).Value = 45218
$ws.Cells.Item(346, 5).Value = 5
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100107
$ws.Cells.Item(346, 8).Value = "Otros"
$ws.Cells.Item(346, 9).Value = 100107002
$ws.Cells.Item(346, 10).Value = "Chirimoya"
$ws.Cells.Item(346, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(346, 12).Value = "Especial"
$ws.Cells.Item(346, 13).Value = 56
$ws.Cells.Item(346, 14).Value = 30000
$ws.Cells.Item(346, 15).Value = 30000
$ws.Cells.Item(346, 16).Value = 30000
$ws.Cells.Item(346, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(346, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(346, 19).Value = 3000
$ws.Cells.Item(346, 20).Value = 10

# Row 347: Primera
$ws.Cells.Item(347, 1).Value = 3
$ws.Cells.Item(347, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(347, 3).Value = "Coquimbo"
$ws.Cells.Item(347, 4).Value = 45218
$ws.Cells.Item(347, 5).Value = 5
$ws.Cells.Item(347, 6).Value = "Fruta"
$ws.Cells.Item(347, 7).Value = 100107
$ws.Cells.Item(347, 8).Value = "Otros"
$ws.Cells.Item(347, 9).Value = 100107002
$ws.Cells.Item(347, 10).Value = "Chirimoya"
$ws.Cells.Item(347, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(347, 12).Value = "Primera"
$ws.Cells.Item(347, 13).Value = 67
$ws.Cells.Item(347, 14).Value = 27000
$ws.Cells.Item(347, 15).Value = 27000
$ws.Cells.Item(347, 16).Value = 27000
$ws.Cells.Item(347, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(347, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(347, 19).Value = 2700
$ws.Cells.Item(347, 20).Value = 10

# Row 348: Segunda
$ws.Cells.Item(348, 1).Value = 3
$ws.Cells.Item(348, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(348, 3).Value = "Coquimbo"
$ws.Cells.Item(348, 4).Value = 45218
$ws.Cells.Item(348, 5).Value = 5
$ws.Cells.Item(348, 6).Value = "Fruta"
$ws.Cells.Item(348, 7).Value = 100107
$ws.Cells.Item(348, 8).Value = "Otros"
$ws.Cells.Item(348, 9).Value = 100107002
$ws.Cells.Item(348, 10).Value = "Chirimoya"
$ws.Cells.Item(348, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(348, 12).Value = "Segunda"
$ws.Cells.Item(348, 13).Value = 56
$ws.Cells.Item(348, 14).Value = 23000
$ws.Cells.Item(348, 15).Value = 23000
$ws.Cells.Item(348, 16).Value = 23000
$ws.Cells.Item(348, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(348, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(348, 19).Value = 2300
$ws.Cells.Item(348, 20).Value = 10
